$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best-effort: reposition/resize the document window to match the saved
# workbookView geometry from the source edit. Real Excel persists this from
# live OS window placement; headless hosts may not round-trip it, so this is
# wrapped defensively and the rest of the script does not depend on it.
try {
    $win = $excel.ActiveWindow
    $win.Left = -51200
    $win.Top = -20240
    $win.Width = 51200
    $win.Height = 28300
} catch {
}

# Update header A1 label ("Code" -> "Code(in Hex)")
$ws.Range("A1").Value = "Code(in Hex)"

# Row 23 was previously an empty gap row (sheet jumped from 22 to 24); just
# populate it directly -- no shifting of row 24 should occur.
$ws.Range("A23").Value = '"015"'
$ws.Range("B23").Value = "Sex"
$ws.Range("C23").Value = 1

# Copy the style used by the preceding block (rows 13-22) onto the new row's code cell.
$ws.Range("A22").Copy()
$ws.Range("A23").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# View settings: zoom to 200% and move the selection to B6.
$excel.ActiveWindow.Zoom = 200
$ws.Range("B6").Select() | Out-Null
